$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update values in column B (rows 2-13) ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 27940.09
$metrics.Range("B3").Value = 22879.260000000002
$metrics.Range("B4").Value = 8809.86
$metrics.Range("B5").Value = 1100
$metrics.Range("B6").Value = 4409181.709999999
$metrics.Range("B7").Value = 3724429.1600000006
$metrics.Range("B8").Value = 1283849.6300000001
$metrics.Range("B9").Value = 170636
$metrics.Range("B10").Value = 32874518.950000003
$metrics.Range("B11").Value = 30999659.949999999
$metrics.Range("B12").Value = 11565563.060000001
$metrics.Range("B13").Value = 1268264

# Update the selection on the Metrics sheet (was E12 -> now D14)
$metrics.Range("D14").Select()

# --- "today" sheet: values are formula-driven off Metrics, so they will
#     recalc automatically; A1 (=TODAY()-1) also recalculates on its own.
#     Only the active selection needs to be moved (was C7 -> now F7).
#     Select this sheet/range last so it remains the active tab, matching
#     the workbook's original active-tab state. ---
$today = $wb.Worksheets.Item("today")
$today.Range("F7").Select()
